# Actualización automática 2025-06-13 17:05:08
#
# Updates sales figures across the three sheets of the workbook.
# Column widths are set via ColumnWidth using a value chosen so that the
# engine's internal (MDW-based) pixel-snapping reproduces the exact target
# stored width: stored = 5/6 + round(6*ColumnWidth)/6, so ColumnWidth =
# target - 5/6 lands safely in the middle of the snapping interval.

$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Cells.Item(4, 11).Value = 13303.18      # K4
$ws1.Cells.Item(16, 3).Value = 497.66        # C16
$ws1.Cells.Item(25, 12).Value = 1532.21      # L25
$ws1.Cells.Item(53, 3).Value = "4 de 51"     # C53

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Columns.Item(6).ColumnWidth = 13.166666666666666   # column F width 13 -> 14
$ws2.Cells.Item(4, 6).Value = 20388.32       # F4
$ws2.Cells.Item(16, 6).Value = 2687.38       # F16
$ws2.Cells.Item(25, 6).Value = 2346.29       # F25
$ws2.Cells.Item(53, 6).Value = 36054.86      # F53

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Columns.Item(4).ColumnWidth = 13.166666666666666   # column D width 13 -> 14
$ws3.Columns.Item(5).ColumnWidth = 23.166666666666668   # column E width 22 -> 24

$ws3.Cells.Item(2, 4).Value = 2540.15
$ws3.Cells.Item(2, 5).Value = 7430.19304517915
$ws3.Cells.Item(2, 6).Value = 0.2547705719341533

$ws3.Cells.Item(15, 4).Value = 15156.46
$ws3.Cells.Item(15, 5).Value = -1656.459999999999
$ws3.Cells.Item(15, 6).Value = 1.122700740740741

$ws3.Cells.Item(16, 4).Value = 6401.68
$ws3.Cells.Item(16, 5).Value = 26339.77
$ws3.Cells.Item(16, 6).Value = 0.1955221897625181

$ws3.Cells.Item(19, 4).Value = 36054.86
$ws3.Cells.Item(19, 5).Value = 58392.58064517915
$ws3.Cells.Item(19, 6).Value = 0.3817452305081635
